# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.581.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.52%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.254.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'302.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.55%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'91.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.86%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.528"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.10%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.47%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'32.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.62%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'52.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.64%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.87%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.97%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +2.77%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.603.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.93%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'14.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.48%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.259.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.36%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D19").Value = "'41.509.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.60%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +8.24%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +1.55%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.10%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'66.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.37%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'239.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.69%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.27%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.01%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +4.72%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'23.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.73%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.73%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.37%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'160.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.66%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'34.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +6.11%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.02%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +3.41%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +3.43%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.12%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.34%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +2.26%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +5.31%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.44%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +2.35%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.04%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.050.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.30%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'19.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.15%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +2.18%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.99%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.24%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +6.13%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +3.44%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'72.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +6.72%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +2.00%  "
$ws.Range("E51").Style = "Normal"
